# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Garuda Profits workbook (updates cached price/profit columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 53
$ws.Range("H53").Value = 41982.5
$ws.Range("I53").Value = 100175.2
$ws.Range("J53").Value = 416.2857
$ws.Range("K53").Value = 100175.2
$ws.Range("L53").Value = 416.2857
$ws.Range("M53").Value = -99538.2
$ws.Range("N53").Value = -1690.2857

# row 94
$ws.Range("H94").Value = 2800.5
$ws.Range("J94").Value = 4000
$ws.Range("L94").Value = 4000
$ws.Range("N94").Value = -4902

# row 107
$ws.Range("H107").Value = 528.73914
$ws.Range("I107").Value = 400.875
$ws.Range("J107").Value = 821
$ws.Range("K107").Value = 400.875
$ws.Range("L107").Value = 821
$ws.Range("M107").Value = 1519.125
$ws.Range("N107").Value = -4661

# row 113
$ws.Range("H113").Value = 1678.5
$ws.Range("I113").Value = 1282.3077
$ws.Range("J113").Value = 2708.6
$ws.Range("K113").Value = 1282.3077
$ws.Range("L113").Value = 2708.6
$ws.Range("M113").Value = 1971.6923
$ws.Range("N113").Value = -9216.6

# row 125
$ws.Range("H125").Value = 1952.6666
$ws.Range("I125").Value = 1706.4
$ws.Range("J125").Value = 2128.5715
$ws.Range("K125").Value = 15357.6
$ws.Range("L125").Value = 19157.1435
$ws.Range("M125").Value = -12897.6
$ws.Range("N125").Value = -24077.1435

# row 131
$ws.Range("H131").Value = 1144.5
$ws.Range("I131").Value = 767.5
$ws.Range("J131").Value = 2652.5
$ws.Range("K131").Value = 2302.5
$ws.Range("L131").Value = 7957.5
$ws.Range("M131").Value = 2737.5
$ws.Range("N131").Value = -18037.5

# row 137
$ws.Range("H137").Value = 1859.6
$ws.Range("I137").Value = 1685.6428
$ws.Range("J137").Value = 2265.5
$ws.Range("K137").Value = 5056.928400000001
$ws.Range("L137").Value = 6796.5
$ws.Range("M137").Value = -2506.928400000001
$ws.Range("N137").Value = -11896.5

# row 141
$ws.Range("H141").Value = 2195.2307
$ws.Range("I141").Value = 1530.5264
$ws.Range("J141").Value = 3999.4285
$ws.Range("K141").Value = 4591.5792
$ws.Range("L141").Value = 11998.2855
$ws.Range("M141").Value = 588.4207999999999
$ws.Range("N141").Value = -22358.2855

$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 2578
$ws.Range("I61").Value = 1156
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1156
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -944
$ws.Range("N61").Value = -4424

# row 74
$ws.Range("H74").Value = 1686.7894
$ws.Range("I74").Value = 1717.7858
$ws.Range("J74").Value = 1600
$ws.Range("K74").Value = 1717.7858
$ws.Range("L74").Value = 1600
$ws.Range("M74").Value = -843.7858000000001
$ws.Range("N74").Value = -3348

# row 77
$ws.Range("H77").Value = 1686.7894
$ws.Range("I77").Value = 1717.7858
$ws.Range("J77").Value = 1600
$ws.Range("K77").Value = 8588.929
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = -4220.929
$ws.Range("N77").Value = -16736

# row 132
$ws.Range("H132").Value = 4618.4375
$ws.Range("I132").Value = 4530.4614
$ws.Range("K132").Value = 13591.3842
$ws.Range("M132").Value = -11061.3842

# row 136
$ws.Range("H136").Value = 2578
$ws.Range("I136").Value = 1156
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3468
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -918
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 39839.5
$ws.Range("I134").Value = 48610.81
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 145832.43
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -143297.43
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 8002832.5
$ws.Range("I31").Value = 2905
$ws.Range("J31").Value = 100002000
$ws.Range("K31").Value = 2905
$ws.Range("L31").Value = 100002000
$ws.Range("M31").Value = -2610
$ws.Range("N31").Value = -100002590

# row 34
$ws.Range("H34").Value = 8002832.5
$ws.Range("I34").Value = 2905
$ws.Range("J34").Value = 100002000
$ws.Range("K34").Value = 2905
$ws.Range("L34").Value = 100002000
$ws.Range("M34").Value = -2703
$ws.Range("N34").Value = -100002404

# row 58
$ws.Range("H58").Value = 1059.75
$ws.Range("I58").Value = 1128
$ws.Range("J58").Value = 946
$ws.Range("K58").Value = 1128
$ws.Range("L58").Value = 946
$ws.Range("M58").Value = -925
$ws.Range("N58").Value = -1352

# row 132
$ws.Range("H132").Value = 1741.6976
$ws.Range("I132").Value = 1176.6
$ws.Range("J132").Value = 3045.7693
$ws.Range("K132").Value = 3529.8
$ws.Range("L132").Value = 9137.3079
$ws.Range("M132").Value = -999.7999999999997
$ws.Range("N132").Value = -14197.3079

# row 134
$ws.Range("H134").Value = 849.9655
$ws.Range("I134").Value = 757.37036
$ws.Range("K134").Value = 2272.11108
$ws.Range("M134").Value = 262.8889199999999

# row 136
$ws.Range("H136").Value = 1059.75
$ws.Range("I136").Value = 1128
$ws.Range("J136").Value = 946
$ws.Range("K136").Value = 3384
$ws.Range("L136").Value = 2838
$ws.Range("M136").Value = -834
$ws.Range("N136").Value = -7938

$ws = $wb.Worksheets.Item("CUL")
# row 97
$ws.Range("H97").Value = 344.33334
$ws.Range("I97").Value = 355.33334
$ws.Range("K97").Value = 1066.00002
$ws.Range("M97").Value = -570.0000199999999

# row 122
$ws.Range("H122").Value = 1828.1111
$ws.Range("J122").Value = 1600.3
$ws.Range("L122").Value = 14402.7
$ws.Range("N122").Value = -19302.7

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 8230387
$ws.Range("I70").Value = 11091146
$ws.Range("J70").Value = 5705.625
$ws.Range("K70").Value = 11091146
$ws.Range("L70").Value = 5705.625
$ws.Range("M70").Value = -11090876
$ws.Range("N70").Value = -6245.625

# row 73
$ws.Range("H73").Value = 8230387
$ws.Range("I73").Value = 11091146
$ws.Range("J73").Value = 5705.625
$ws.Range("K73").Value = 11091146
$ws.Range("L73").Value = 5705.625
$ws.Range("M73").Value = -11090210
$ws.Range("N73").Value = -7577.625

# row 107
$ws.Range("H107").Value = 589.2414
$ws.Range("I107").Value = 490.2381
$ws.Range("J107").Value = 849.125
$ws.Range("K107").Value = 490.2381
$ws.Range("L107").Value = 849.125
$ws.Range("M107").Value = 1429.7619
$ws.Range("N107").Value = -4689.125

# row 132
$ws.Range("H132").Value = 103041.1
$ws.Range("I132").Value = 288432
$ws.Range("J132").Value = 3215.2307
$ws.Range("K132").Value = 865296
$ws.Range("L132").Value = 9645.6921
$ws.Range("M132").Value = -862766
$ws.Range("N132").Value = -14705.6921

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 3242.3076
$ws.Range("I7").Value = 2975.8333
$ws.Range("J7").Value = 3470.7144
$ws.Range("K7").Value = 2975.8333
$ws.Range("L7").Value = 3470.7144
$ws.Range("M7").Value = -2863.8333
$ws.Range("N7").Value = -3694.7144

# row 126
$ws.Range("H126").Value = 3242.3076
$ws.Range("I126").Value = 2975.8333
$ws.Range("J126").Value = 3470.7144
$ws.Range("K126").Value = 8927.499899999999
$ws.Range("L126").Value = 10412.1432
$ws.Range("M126").Value = -6457.499899999999
$ws.Range("N126").Value = -15352.1432

$ws = $wb.Worksheets.Item("WVR")
# row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# row 94
$ws.Range("H94").Value = 52400
$ws.Range("J94").Value = 52400
$ws.Range("L94").Value = 52400
$ws.Range("N94").Value = -54202

# row 126
$ws.Range("H126").Value = 1944.909
$ws.Range("I126").Value = 2598.5
$ws.Range("J126").Value = 1571.4286
$ws.Range("K126").Value = 7795.5
$ws.Range("L126").Value = 4714.2858
$ws.Range("M126").Value = -5325.5
$ws.Range("N126").Value = -9654.2858

# row 132
$ws.Range("H132").Value = 1044.5483
$ws.Range("I132").Value = 771.88
$ws.Range("J132").Value = 2180.6667
$ws.Range("K132").Value = 2315.64
$ws.Range("L132").Value = 6542.000100000001
$ws.Range("M132").Value = 214.3600000000001
$ws.Range("N132").Value = -11602.0001

# row 136
$ws.Range("H136").Value = 6562.55
$ws.Range("I136").Value = 7561.8237
$ws.Range("K136").Value = 22685.4711
$ws.Range("M136").Value = -20135.4711
